# Updated excel thermister sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename the "Voltage Divider Low" header (now in K1) to "Voltage Divider"
$ws.Range("K1").Value2 = "Voltage Divider"

# Re-create the voltage-divider formulas in column K (was column M)
for ($r = 2; $r -le 17; $r++) {
  $ws.Range("K$r").Formula = $ws.Range("M$r").Formula
}

# Clear the old column M, the data now lives in column K
$ws.Range("M1:M17").Clear()

# Update the active selection to match the saved view
$ws.Range("U14").Select()
